# Apply "write reporting log files" edit to the active workbook.
#
# Rule observed in the data:
#   - For rows where "Did Harvest Occur?" (column B) is "No" (and Species/
#     column F is blank, Unknown Sex Count/column J is 0), flip them to
#     reflect that a harvest log was actually written:
#       B -> "Yes"
#       F -> "Na"
#       J -> 1
#   - For rows that already recorded a species in column F (all caps),
#     normalize the text to Title Case (e.g. "WEASEL" -> "Weasel").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)   # column B - Did Harvest Occur?
    $fCell = $ws.Cells.Item($r, 6)   # column F - Species
    $jCell = $ws.Cells.Item($r, 10)  # column J - Unknown Sex Count

    $bVal = $bCell.Value()

    if ($bVal -eq "No") {
        $bCell.Value = "Yes"
        $fCell.Value = "Na"
        $jCell.Value = 1
    }
    else {
        $fVal = $fCell.Value()
        if ($fVal -ne $null -and $fVal -ne "") {
            $titleCased = $fVal.Substring(0,1).ToUpper() + $fVal.Substring(1).ToLower()
            if (-not $titleCased.Equals($fVal)) {
                $fCell.Value = $titleCased
            }
        }
    }
}
